# chore: update Sheets via scheduled runner
#
# Refreshes computed profit/loss figures (columns H-N) across several rows
# on each of the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 6709.1562  # H19: was 6581.9697
$ws.Cells.Item(19, 9).Value = 472.8  # I19: was 472.9
$ws.Cells.Item(19, 10).Value = 9543.862999999999  # J19: was 9238.087
$ws.Cells.Item(19, 11).Value = 472.8  # K19: was 472.9
$ws.Cells.Item(19, 12).Value = 9543.862999999999  # L19: was 9238.087
$ws.Cells.Item(19, 13).Value = -297.8  # M19: was -297.9
$ws.Cells.Item(19, 14).Value = -9893.862999999999  # N19: was -9588.087
$ws.Cells.Item(21, 8).Value = 17285.834  # H21: was 17143.2
$ws.Cells.Item(21, 10).Value = 18539.6  # J21: was 18674.75
$ws.Cells.Item(21, 12).Value = 18539.6  # L21: was 18674.75
$ws.Cells.Item(21, 14).Value = -19475.6  # N21: was -19610.75
$ws.Cells.Item(23, 8).Value = 17285.834  # H23: was 17143.2
$ws.Cells.Item(23, 10).Value = 18539.6  # J23: was 18674.75
$ws.Cells.Item(23, 12).Value = 18539.6  # L23: was 18674.75
$ws.Cells.Item(23, 14).Value = -19007.6  # N23: was -19142.75
$ws.Cells.Item(38, 8).Value = 477.94116  # H38: was 557.8125
$ws.Cells.Item(38, 9).Value = 387.5  # I38: was 443.75
$ws.Cells.Item(38, 11).Value = 1162.5  # K38: was 1331.25
$ws.Cells.Item(38, 13).Value = -790.5  # M38: was -959.25
$ws.Cells.Item(58, 8).Value = 1195.1818  # H58: was 1179.4762
$ws.Cells.Item(58, 9).Value = 846.26666  # I58: was 857.93335
$ws.Cells.Item(58, 10).Value = 1942.8572  # J58: was 1983.3334
$ws.Cells.Item(58, 11).Value = 2538.79998  # K58: was 2573.80005
$ws.Cells.Item(58, 12).Value = 5828.571599999999  # L58: was 5950.0002
$ws.Cells.Item(58, 13).Value = -2388.79998  # M58: was -2423.80005
$ws.Cells.Item(58, 14).Value = -6128.571599999999  # N58: was -6250.0002
$ws.Cells.Item(70, 8).Value = 998.8333  # H70: was 998.25
$ws.Cells.Item(70, 10).Value = 998.25  # J70: was 996.5
$ws.Cells.Item(70, 12).Value = 2994.75  # L70: was 2989.5
$ws.Cells.Item(70, 14).Value = -3534.75  # N70: was -3529.5
$ws.Cells.Item(73, 8).Value = 998.8333  # H73: was 998.25
$ws.Cells.Item(73, 10).Value = 998.25  # J73: was 996.5
$ws.Cells.Item(73, 12).Value = 2994.75  # L73: was 2989.5
$ws.Cells.Item(73, 14).Value = -4866.75  # N73: was -4861.5
$ws.Cells.Item(112, 8).Value = 8248.68  # H112: was 8035.269
$ws.Cells.Item(112, 10).Value = 8248.68  # J112: was 8035.269
$ws.Cells.Item(112, 12).Value = 24746.04  # L112: was 24105.807
$ws.Cells.Item(112, 14).Value = -26962.04  # N112: was -26321.807
$ws.Cells.Item(121, 8).Value = 395  # H121: was 2590
$ws.Cells.Item(121, 9).Value = 395  # I121: was 380
$ws.Cells.Item(121, 10).Value = 0  # J121: was 4800
$ws.Cells.Item(121, 11).Value = 1185  # K121: was 1140
$ws.Cells.Item(121, 12).Value = 0  # L121: was 14400
$ws.Cells.Item(121, 13).Value = 562  # M121: was 607
$ws.Cells.Item(121, 14).ClearContents()  # N121: was -17894
$ws.Cells.Item(129, 8).Value = 1001.32465  # H129: was 1023.1739
$ws.Cells.Item(129, 9).Value = 542.61536  # I129: was 571.1667
$ws.Cells.Item(129, 10).Value = 1094.5  # J129: was 1090.975
$ws.Cells.Item(129, 11).Value = 1627.84608  # K129: was 1713.5001
$ws.Cells.Item(129, 12).Value = 3283.5  # L129: was 3272.925
$ws.Cells.Item(129, 13).Value = 3372.15392  # M129: was 3286.4999
$ws.Cells.Item(129, 14).Value = -13283.5  # N129: was -13272.925
$ws.Cells.Item(137, 8).Value = 3642.077  # H137: was 4213.8184
$ws.Cells.Item(137, 9).Value = 2354.3044  # I137: was 2745.2104
$ws.Cells.Item(137, 11).Value = 7062.9132  # K137: was 8235.6312
$ws.Cells.Item(137, 13).Value = -4512.9132  # M137: was -5685.6312
$ws.Cells.Item(138, 8).Value = 1909.2727  # H138: was 2016.4
$ws.Cells.Item(138, 10).Value = 2108.5264  # J138: was 2447.3333
$ws.Cells.Item(138, 12).Value = 6325.5792  # L138: was 7341.999899999999
$ws.Cells.Item(138, 14).Value = -16605.5792  # N138: was -17621.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3517.2  # H61: was 3736.05
$ws.Cells.Item(61, 9).Value = 2800.4285  # I61: was 2946.7778
$ws.Cells.Item(61, 10).Value = 4144.375  # J61: was 4381.8184
$ws.Cells.Item(61, 11).Value = 2800.4285  # K61: was 2946.7778
$ws.Cells.Item(61, 12).Value = 4144.375  # L61: was 4381.8184
$ws.Cells.Item(61, 13).Value = -2588.4285  # M61: was -2734.7778
$ws.Cells.Item(61, 14).Value = -4568.375  # N61: was -4805.8184
$ws.Cells.Item(74, 8).Value = 1373.9375  # H74: was 1405.3125
$ws.Cells.Item(74, 9).Value = 968  # I74: was 1052.7693
$ws.Cells.Item(74, 10).Value = 3133  # J74: was 2933
$ws.Cells.Item(74, 11).Value = 968  # K74: was 1052.7693
$ws.Cells.Item(74, 12).Value = 3133  # L74: was 2933
$ws.Cells.Item(74, 13).Value = -94  # M74: was -178.7692999999999
$ws.Cells.Item(74, 14).Value = -4881  # N74: was -4681
$ws.Cells.Item(77, 8).Value = 1373.9375  # H77: was 1405.3125
$ws.Cells.Item(77, 9).Value = 968  # I77: was 1052.7693
$ws.Cells.Item(77, 10).Value = 3133  # J77: was 2933
$ws.Cells.Item(77, 11).Value = 4840  # K77: was 5263.8465
$ws.Cells.Item(77, 12).Value = 15665  # L77: was 14665
$ws.Cells.Item(77, 13).Value = -472  # M77: was -895.8464999999997
$ws.Cells.Item(77, 14).Value = -24401  # N77: was -23401
$ws.Cells.Item(132, 8).Value = 3622.919  # H132: was 3978.6667
$ws.Cells.Item(132, 9).Value = 2193.6924  # I132: was 2467.4546
$ws.Cells.Item(132, 11).Value = 6581.0772  # K132: was 7402.3638
$ws.Cells.Item(132, 13).Value = -4051.0772  # M132: was -4872.3638
$ws.Cells.Item(136, 8).Value = 3517.2  # H136: was 3736.05
$ws.Cells.Item(136, 9).Value = 2800.4285  # I136: was 2946.7778
$ws.Cells.Item(136, 10).Value = 4144.375  # J136: was 4381.8184
$ws.Cells.Item(136, 11).Value = 8401.2855  # K136: was 8840.3334
$ws.Cells.Item(136, 12).Value = 12433.125  # L136: was 13145.4552
$ws.Cells.Item(136, 13).Value = -5851.2855  # M136: was -6290.3334
$ws.Cells.Item(136, 14).Value = -17533.125  # N136: was -18245.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2656.2334  # H134: was 2777.8965
$ws.Cells.Item(134, 9).Value = 2535.762  # I134: was 2706.15
$ws.Cells.Item(134, 11).Value = 7607.286  # K134: was 8118.450000000001
$ws.Cells.Item(134, 13).Value = -5072.286  # M134: was -5583.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 28900.2  # H4: was 28091.092
$ws.Cells.Item(4, 10).Value = 28900.2  # J4: was 28091.092
$ws.Cells.Item(4, 12).Value = 28900.2  # L4: was 28091.092
$ws.Cells.Item(4, 14).Value = -29124.2  # N4: was -28315.092
$ws.Cells.Item(31, 8).Value = 5591.354  # H31: was 5146.9575
$ws.Cells.Item(31, 9).Value = 1000.43243  # I31: was 979.7368
$ws.Cells.Item(31, 10).Value = 11657.929  # J31: was 9945.575999999999
$ws.Cells.Item(31, 11).Value = 1000.43243  # K31: was 979.7368
$ws.Cells.Item(31, 12).Value = 11657.929  # L31: was 9945.575999999999
$ws.Cells.Item(31, 13).Value = -705.43243  # M31: was -684.7368
$ws.Cells.Item(31, 14).Value = -12247.929  # N31: was -10535.576
$ws.Cells.Item(34, 8).Value = 5591.354  # H34: was 5146.9575
$ws.Cells.Item(34, 9).Value = 1000.43243  # I34: was 979.7368
$ws.Cells.Item(34, 10).Value = 11657.929  # J34: was 9945.575999999999
$ws.Cells.Item(34, 11).Value = 1000.43243  # K34: was 979.7368
$ws.Cells.Item(34, 12).Value = 11657.929  # L34: was 9945.575999999999
$ws.Cells.Item(34, 13).Value = -798.43243  # M34: was -777.7368
$ws.Cells.Item(34, 14).Value = -12061.929  # N34: was -10349.576
$ws.Cells.Item(58, 8).Value = 1843.1538  # H58: was 1676.8667
$ws.Cells.Item(58, 9).Value = 1683.1111  # I58: was 1436.5834
$ws.Cells.Item(58, 10).Value = 2203.25  # J58: was 2638
$ws.Cells.Item(58, 11).Value = 1683.1111  # K58: was 1436.5834
$ws.Cells.Item(58, 12).Value = 2203.25  # L58: was 2638
$ws.Cells.Item(58, 13).Value = -1480.1111  # M58: was -1233.5834
$ws.Cells.Item(58, 14).Value = -2609.25  # N58: was -3044
$ws.Cells.Item(132, 8).Value = 5378062  # H132: was 6946514
$ws.Cells.Item(132, 9).Value = 1471.2307  # I132: was 1758.55
$ws.Cells.Item(132, 10).Value = 33336334  # J132: was 41670292
$ws.Cells.Item(132, 11).Value = 4413.6921  # K132: was 5275.65
$ws.Cells.Item(132, 12).Value = 100009002  # L132: was 125010876
$ws.Cells.Item(132, 13).Value = -1883.6921  # M132: was -2745.65
$ws.Cells.Item(132, 14).Value = -100014062  # N132: was -125015936
$ws.Cells.Item(134, 8).Value = 2125.6538  # H134: was 1858.2903
$ws.Cells.Item(134, 9).Value = 1243.7273  # I134: was 1203.2174
$ws.Cells.Item(134, 10).Value = 6976.25  # J134: was 3741.625
$ws.Cells.Item(134, 11).Value = 3731.1819  # K134: was 3609.6522
$ws.Cells.Item(134, 12).Value = 20928.75  # L134: was 11224.875
$ws.Cells.Item(134, 13).Value = -1196.1819  # M134: was -1074.6522
$ws.Cells.Item(134, 14).Value = -25998.75  # N134: was -16294.875
$ws.Cells.Item(136, 8).Value = 1843.1538  # H136: was 1676.8667
$ws.Cells.Item(136, 9).Value = 1683.1111  # I136: was 1436.5834
$ws.Cells.Item(136, 10).Value = 2203.25  # J136: was 2638
$ws.Cells.Item(136, 11).Value = 5049.3333  # K136: was 4309.7502
$ws.Cells.Item(136, 12).Value = 6609.75  # L136: was 7914
$ws.Cells.Item(136, 13).Value = -2499.3333  # M136: was -1759.7502
$ws.Cells.Item(136, 14).Value = -11709.75  # N136: was -13014

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 322.85715  # H5: was 335.76923
$ws.Cells.Item(5, 9).Value = 297.77777  # I5: was 309.2
$ws.Cells.Item(5, 11).Value = 893.33331  # K5: was 927.5999999999999
$ws.Cells.Item(5, 13).Value = -781.33331  # M5: was -815.5999999999999
$ws.Cells.Item(12, 9).Value = 5  # I12: was 4
$ws.Cells.Item(12, 10).Value = 83.888885  # J12: was 94
$ws.Cells.Item(12, 11).Value = 15  # K12: was 12
$ws.Cells.Item(12, 12).Value = 251.666655  # L12: was 282
$ws.Cells.Item(12, 13).Value = 158  # M12: was 161
$ws.Cells.Item(12, 14).Value = -597.666655  # N12: was -628
$ws.Cells.Item(110, 8).Value = 13492  # H110: was 13496.348
$ws.Cells.Item(110, 10).Value = 13568.954  # J110: was 13573.5
$ws.Cells.Item(110, 12).Value = 40706.862  # L110: was 40720.5
$ws.Cells.Item(110, 14).Value = -48886.862  # N110: was -48900.5
$ws.Cells.Item(122, 8).Value = 5579.4  # H122: was 5597.1
$ws.Cells.Item(122, 9).Value = 529.61536  # I122: was 545.6429000000001
$ws.Cells.Item(122, 10).Value = 14957.571  # J122: was 17383.834
$ws.Cells.Item(122, 11).Value = 4766.53824  # K122: was 4910.7861
$ws.Cells.Item(122, 12).Value = 134618.139  # L122: was 156454.506
$ws.Cells.Item(122, 13).Value = -2316.53824  # M122: was -2460.7861
$ws.Cells.Item(122, 14).Value = -139518.139  # N122: was -161354.506
$ws.Cells.Item(131, 8).Value = 1056.6721  # H131: was 1080.7954
$ws.Cells.Item(131, 10).Value = 1061.4746  # J131: was 1088.6904
$ws.Cells.Item(131, 12).Value = 3184.4238  # L131: was 3266.0712
$ws.Cells.Item(131, 14).Value = -13264.4238  # N131: was -13346.0712
$ws.Cells.Item(135, 8).Value = 322.85715  # H135: was 335.76923
$ws.Cells.Item(135, 9).Value = 297.77777  # I135: was 309.2
$ws.Cells.Item(135, 11).Value = 2679.99993  # K135: was 2782.8
$ws.Cells.Item(135, 13).Value = -144.9999299999999  # M135: was -247.7999999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 7514999  # H57: was 10012666
$ws.Cells.Item(57, 10).Value = 7514999  # J57: was 10012666
$ws.Cells.Item(57, 12).Value = 7514999  # L57: was 10012666
$ws.Cells.Item(57, 14).Value = -7516639  # N57: was -10014306
$ws.Cells.Item(80, 8).Value = 50942804  # H80: was 50942904
$ws.Cells.Item(80, 9).Value = 63628130  # I80: was 72717576
$ws.Cells.Item(80, 10).Value = 201500  # J80: was 135333.33
$ws.Cells.Item(80, 11).Value = 63628130  # K80: was 72717576
$ws.Cells.Item(80, 12).Value = 201500  # L80: was 135333.33
$ws.Cells.Item(80, 13).Value = -63627132  # M80: was -72716578
$ws.Cells.Item(80, 14).Value = -203496  # N80: was -137329.33
$ws.Cells.Item(83, 8).Value = 50942804  # H83: was 50942904
$ws.Cells.Item(83, 9).Value = 63628130  # I83: was 72717576
$ws.Cells.Item(83, 10).Value = 201500  # J83: was 135333.33
$ws.Cells.Item(83, 11).Value = 318140650  # K83: was 363587880
$ws.Cells.Item(83, 12).Value = 1007500  # L83: was 676666.6499999999
$ws.Cells.Item(83, 13).Value = -318135658  # M83: was -363582888
$ws.Cells.Item(83, 14).Value = -1017484  # N83: was -686650.6499999999
$ws.Cells.Item(132, 8).Value = 3421.2083  # H132: was 3275.5
$ws.Cells.Item(132, 9).Value = 3500.8572  # I132: was 3100.875
$ws.Cells.Item(132, 10).Value = 3309.7  # J132: was 3624.75
$ws.Cells.Item(132, 11).Value = 10502.5716  # K132: was 9302.625
$ws.Cells.Item(132, 12).Value = 9929.099999999999  # L132: was 10874.25
$ws.Cells.Item(132, 13).Value = -7972.571599999999  # M132: was -6772.625
$ws.Cells.Item(132, 14).Value = -14989.1  # N132: was -15934.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 2048000.8  # H2: was 1713334
$ws.Cells.Item(2, 10).Value = 2048000.8  # J2: was 1713334
$ws.Cells.Item(2, 12).Value = 2048000.8  # L2: was 1713334
$ws.Cells.Item(2, 14).Value = -2048224.8  # N2: was -1713558

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 61249.75  # H64: was 0
$ws.Cells.Item(64, 10).Value = 61249.75  # J64: was 0
$ws.Cells.Item(64, 12).Value = 61249.75  # L64: was 0
$ws.Cells.Item(64, 14).Value = -61745.75  # N64: was None
$ws.Cells.Item(67, 8).Value = 61249.75  # H67: was 0
$ws.Cells.Item(67, 10).Value = 61249.75  # J67: was 0
$ws.Cells.Item(67, 12).Value = 61249.75  # L67: was 0
$ws.Cells.Item(67, 14).Value = -62965.75  # N67: was None
$ws.Cells.Item(81, 8).Value = 6289  # H81: was 6385.4443
$ws.Cells.Item(81, 9).Value = 5209.857  # I81: was 5578.1665
$ws.Cells.Item(81, 10).Value = 7799.8  # J81: was 8000
$ws.Cells.Item(81, 11).Value = 10419.714  # K81: was 11156.333
$ws.Cells.Item(81, 12).Value = 15599.6  # L81: was 16000
$ws.Cells.Item(81, 13).Value = -9358.714  # M81: was -10095.333
$ws.Cells.Item(81, 14).Value = -17721.6  # N81: was -18122
$ws.Cells.Item(84, 8).Value = 6289  # H84: was 6385.4443
$ws.Cells.Item(84, 9).Value = 5209.857  # I84: was 5578.1665
$ws.Cells.Item(84, 10).Value = 7799.8  # J84: was 8000
$ws.Cells.Item(84, 11).Value = 52098.57  # K84: was 55781.665
$ws.Cells.Item(84, 12).Value = 77998  # L84: was 80000
$ws.Cells.Item(84, 13).Value = -46794.57  # M84: was -50477.665
$ws.Cells.Item(84, 14).Value = -88606  # N84: was -90608
$ws.Cells.Item(132, 8).Value = 4904347.5  # H132: was 4904326.5
$ws.Cells.Item(132, 9).Value = 2721.389  # I132: was 2684.2632
$ws.Cells.Item(132, 10).Value = 10418677  # J132: was 11113073
$ws.Cells.Item(132, 11).Value = 8164.167  # K132: was 8052.7896
$ws.Cells.Item(132, 12).Value = 31256031  # L132: was 33339219
$ws.Cells.Item(132, 13).Value = -5634.167  # M132: was -5522.7896
$ws.Cells.Item(132, 14).Value = -31261091  # N132: was -33344279
$ws.Cells.Item(136, 8).Value = 2264.3403  # H136: was 2255.96
$ws.Cells.Item(136, 9).Value = 1998.6774  # I136: was 1925.3611
$ws.Cells.Item(136, 10).Value = 2779.0625  # J136: was 3106.0715
$ws.Cells.Item(136, 11).Value = 5996.0322  # K136: was 5776.0833
$ws.Cells.Item(136, 12).Value = 8337.1875  # L136: was 9318.2145
$ws.Cells.Item(136, 13).Value = -3446.0322  # M136: was -3226.0833
$ws.Cells.Item(136, 14).Value = -13437.1875  # N136: was -14418.2145
